$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text representation (these look like
# numbers/dates to Excels auto-detection, e.g. "1.002" or "31.079.83").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '31.079.83'
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").Value = '1.961.16'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '245.25'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").Value = '0.4888'
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("D8").Value = '0.2960'
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("D9").Value = '0.06883'
$ws.Range("E9").Value = '  +1.82%  '
$ws.Range("D10").Value = '19.39'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").Value = '107.18'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").Value = '1.967.06'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '0.07790'
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("D14").Value = '5.466'
$ws.Range("E14").Value = '  +1.06%  '
$ws.Range("D15").Value = '0.7007'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '282.39'
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").Value = '31.098.54'
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").Value = '2.244.33'
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D19").Value = '13.24'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '0.000007730'
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("D22").Value = '5.509'
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '6.517'
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("D25").Value = '9.835'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '168.43'
$ws.Range("E26").Value = '  -1.73%  '
$ws.Range("D27").Value = '20.01'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '2.198'
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("D29").Value = '0.1054'
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("D30").Value = '1.394'
$ws.Range("E30").Value = '  -2.89%  '
$ws.Range("D31").Value = '1.581'
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").Value = '4.619'
$ws.Range("E32").Value = '  -3.29%  '
$ws.Range("D33").Value = '4.447'
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("D34").Value = '0.04946'
$ws.Range("E34").Value = '  -2.49%  '
$ws.Range("D35").Value = '0.7563'
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("D36").Value = '1.171'
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").Value = '2.736'
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("D38").Value = '0.02015'
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("D39").Value = '2.706'
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").Value = '6.535'
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("D41").Value = '78.22'
$ws.Range("E41").Value = '  +12.00%  '
$ws.Range("D42").Value = '2.124'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("D46").Value = '8.132'
$ws.Range("E46").Value = '  +8.79%  '
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '1.031.17'
$ws.Range("E48").Value = '  +11.14%  '
$ws.Range("D51").Value = '35.96'
$ws.Range("E51").Value = '  +0.55%  '

# Rows 44/45 swap places (TheSandbox <-> Quant)
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '109.63'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4467'
$ws.Range("E45").Value = '  +0.35%  '

# Rows 49/50 swap places (Algorand <-> EnergySwap)
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.435'
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1261'
$ws.Range("E50").Value = '  -0.40%  '
